# Update "by_specialty" sheet with new survey data ("eipes"):
# - one new specialty row is inserted (alphabetically) between
#   "ΠΑΙΔΙΑΤΡΙΚΗ" and "ΠΝΕΥΜΟΝΟΛΟΓΙΑ - ΦΥΜΑΤΙΟΛΟΓΙΑ"
# - every row's "n" and "percentage" values are refreshed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23 (shifts existing rows 23-30 down to 24-31)
$ws.Rows.Item(23).Insert()

# New specialty row
$ws.Cells.Item(23, 1).Value = "ΠΛΑΣΤΙΚΗ, ΕΠΑΝΟΡΘΩΤΙΚΗ ΚΑΙ ΑΙΣΘΗΤΙΚΗ ΧΕΙΡΟΥΡΓΙΚΗ"
$ws.Cells.Item(23, 2).Value = 1
$ws.Cells.Item(23, 3).Value = 0

# Refresh n / percentage for every data row (2-31) to match the new dataset
$data = @(
  @(2, 2, 0.01),
  @(3, 21, 0.07),
  @(4, 1, 0),
  @(5, 6, 0.02),
  @(6, 9, 0.03),
  @(7, 10, 0.03),
  @(8, 15, 0.05),
  @(9, 2, 0.01),
  @(10, 7, 0.02),
  @(11, 5, 0.02),
  @(12, 2, 0.01),
  @(13, 21, 0.07),
  @(14, 8, 0.03),
  @(15, 26, 0.08),
  @(16, 8, 0.03),
  @(17, 13, 0.04),
  @(18, 6, 0.02),
  @(19, 9, 0.03),
  @(20, 7, 0.02),
  @(21, 2, 0.01),
  @(22, 23, 0.07),
  @(24, 5, 0.02),
  @(25, 4, 0.01),
  @(26, 16, 0.05),
  @(27, 1, 0),
  @(28, 2, 0.01),
  @(29, 65, 0.2),
  @(30, 14, 0.04),
  @(31, 7, 0.02)
)

foreach ($row in $data) {
  $r = $row[0]
  $n = $row[1]
  $pct = $row[2]
  $ws.Cells.Item($r, 2).Value = $n
  $ws.Cells.Item($r, 3).Value = $pct
}
